$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Cells.Item(771, 1).Value = 'cancellation._$dead'
$ws.Cells.Item(772, 1).Value = 'cancellation._$moved'
$ws.Cells.Item(773, 1).Value = 'cancellation._$circumstances'
$ws.Cells.Item(774, 1).Value = 'cancellation._$abandoned'
$ws.Cells.Item(775, 1).Value = 'cancellation._$unwilling'
$ws.Cells.Item(776, 1).Value = 'risk.health.requirement._$Malaria'
$ws.Cells.Item(777, 1).Value = 'risk.health.requirement._$Medical'
$ws.Cells.Item(778, 1).Value = 'risk.health.requirement._$Wound'
$ws.Cells.Item(779, 1).Value = 'risk.health.requirement._$Mobility'
$ws.Cells.Item(780, 1).Value = 'risk.health.requirement._$Physiotherapy'
$ws.Cells.Item(781, 1).Value = 'risk.health.goal._$Medical'
$ws.Cells.Item(782, 1).Value = 'risk.health.goal._$Sores'
$ws.Cells.Item(783, 1).Value = 'risk.health.goal._$Mobility'
$ws.Cells.Item(784, 1).Value = 'risk.health.goal._$Pain'
$ws.Cells.Item(785, 1).Value = 'risk.social.requirement._$Contact'
$ws.Cells.Item(786, 1).Value = 'risk.social.requirement._$Inclusion'
$ws.Cells.Item(787, 1).Value = 'risk.social.requirement._$Training'
$ws.Cells.Item(788, 1).Value = 'risk.social.goal._$Community'
$ws.Cells.Item(789, 1).Value = 'risk.social.goal._$Inclusion'
$ws.Cells.Item(790, 1).Value = 'risk.nutrition.requirement._$Diabetic'
$ws.Cells.Item(791, 1).Value = 'risk.nutrition.requirement._$Allergies'
$ws.Cells.Item(792, 1).Value = 'risk.nutrition.requirement._$Malnutrition'
$ws.Cells.Item(793, 1).Value = 'risk.nutrition.requirement._$Training'
$ws.Cells.Item(794, 1).Value = 'risk.nutrition.requirement._$Agricultural'
$ws.Cells.Item(795, 1).Value = 'risk.nutrition.goal._$ManageDiet'
$ws.Cells.Item(796, 1).Value = 'risk.nutrition.goal._$Sustainable'
$ws.Cells.Item(797, 1).Value = 'risk.education.requirement._$School'
$ws.Cells.Item(798, 1).Value = 'risk.education.requirement._$Vocational'
$ws.Cells.Item(799, 1).Value = 'risk.education.requirement._$Family'
$ws.Cells.Item(800, 1).Value = 'risk.education.goal._$Education'
$ws.Cells.Item(801, 1).Value = 'risk.education.goal._$Income'
$ws.Cells.Item(802, 1).Value = 'risk.education.goal._$Child'
$ws.Cells.Item(803, 1).Value = 'risk.mental.requirement._$Medical'
$ws.Cells.Item(804, 1).Value = 'risk.mental.requirement._$Family'
$ws.Cells.Item(805, 1).Value = 'risk.mental.requirement._$Community'
$ws.Cells.Item(806, 1).Value = 'risk.mental.requirement._$Refer'
$ws.Cells.Item(807, 1).Value = 'risk.mental.requirement._$CounsellingStressTrauma'
$ws.Cells.Item(808, 1).Value = 'risk.mental.requirement._$CounsellingOther'
$ws.Cells.Item(809, 1).Value = 'risk.mental.goal._$Family'
$ws.Cells.Item(810, 1).Value = 'risk.mental.goal._$Community'
$ws.Cells.Item(811, 1).Value = 'risk.mental.goal._$Medical'
$ws.Cells.Item(812, 1).Value = 'risk.mental.goal._$Agency'
$ws.Cells.Item(813, 1).Value = 'risk.mental.goal._$Improvement'
$ws.Cells.Item(814, 1).Value = 'risk.mental.goal._$Reduction'
$ws.Cells.Item(815, 1).Value = 'risk.mental.goal._$Improvements'
